$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34:F34").Copy()
$ws.Range("A35:F35").PasteSpecial(-4122)

$ws.Range("A35").Value = "wu, wu"
$ws.Range("B35").Value = 2017
$ws.Range("C35").Value = "think globally, fit locally under the manifold setup. Asymptotic analysis of lle"
$ws.Range("D35").Value = "paper"
$ws.Range("E35").Value = "lle, laplace-beltrami"
$ws.Range("F35").Value = "lle not related to laplace-beltrami?!"

$excel.ActiveWindow.ScrollRow = 16
$ws.Range("F23").Select()
